$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has two headers (default + first page) and two footers
# (default + first page), each carrying a single inline logo image.
# The commit renames the logo images' display name:
#   - the Pearson logo (in both footers) from "image2.png" to "image1.png"
#   - the BTEC logo   (in both headers) from "image1.jpg" to "image2.jpg"
# (the underlying embedded picture / relationship target is untouched --
# only the drawing's display name changes). InlineShape.Name does not
# read back its current value, so the BTEC/Pearson logos are told apart
# via their (readable) AlternativeText instead.

for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }

    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
